# Insert 4 new rows at 231-234 (pushes existing rows 231-320 down to 235-324)
# and populate them with a new weekly price record (date 45119 / 2023-07-12)
# for "Agrícola del Norte S.A. de Arica" - Piña - Caramelo - Ecuador.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("231:234").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad  = "Caramelo"
$origen    = "Ecuador"
$fecha     = 45119

# Row 231 - Especial
$r = 231
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = $fecha
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Especial"
$ws.Cells.Item($r, 13).Value2 = 200
$ws.Cells.Item($r, 14).Value2 = 22000
$ws.Cells.Item($r, 15).Value2 = 23000
$ws.Cells.Item($r, 16).Value2 = 22500
$ws.Cells.Item($r, 17).Value2 = "$/caja 10 unidades"
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 2250
$ws.Cells.Item($r, 20).Value2 = 10

# Row 232 - Primera
$r = 232
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = $fecha
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Primera"
$ws.Cells.Item($r, 13).Value2 = 250
$ws.Cells.Item($r, 14).Value2 = 22000
$ws.Cells.Item($r, 15).Value2 = 23000
$ws.Cells.Item($r, 16).Value2 = 22500
$ws.Cells.Item($r, 17).Value2 = "$/caja 12 unidades"
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 1875
$ws.Cells.Item($r, 20).Value2 = 12

# Row 233 - Segunda
$r = 233
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = $fecha
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 22000
$ws.Cells.Item($r, 15).Value2 = 23000
$ws.Cells.Item($r, 16).Value2 = 22500
$ws.Cells.Item($r, 17).Value2 = "$/caja 14 unidades"
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 1607
$ws.Cells.Item($r, 20).Value2 = 14

# Row 234 - Tercera
$r = 234
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = $fecha
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $tipo
$ws.Cells.Item($r, 7).Value2  = $productoId
$ws.Cells.Item($r, 8).Value2  = $producto
$ws.Cells.Item($r, 9).Value2  = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Tercera"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 22000
$ws.Cells.Item($r, 15).Value2 = 23000
$ws.Cells.Item($r, 16).Value2 = 22500
$ws.Cells.Item($r, 17).Value2 = "$/caja 16 unidades"
$ws.Cells.Item($r, 18).Value2 = $origen
$ws.Cells.Item($r, 19).Value2 = 1406
$ws.Cells.Item($r, 20).Value2 = 16
